# ------------------------------------------------------------------
# Apply "[ADDITIONAL SCRAPING]" edit:
#   1. Insert a brand-new "Player Info" worksheet as the FIRST sheet
#      (in front of the existing "ODI Batting" sheet) containing basic
#      player metadata.
#   2. On the existing "ODI Batting" sheet, replace the
#      MATCH_CARD_LINK column (full scorecard URL) with a MATCH_CODE
#      column that only keeps the numeric match code.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- 1. Create the new "Player Info" sheet, inserted before the
#        currently active ("ODI Batting") sheet so it becomes sheet #1 ---
$playerSheet = $wb.Worksheets.Add()
$playerSheet.Name = "Player Info"

# Header row
$playerSheet.Cells.Item(1, 1).Value = "ID"
$playerSheet.Cells.Item(1, 2).Value = "NAME"
$playerSheet.Cells.Item(1, 3).Value = "BATTING_HAND"
$playerSheet.Cells.Item(1, 4).Value = "BOWL_STYLE"

# Style the header row the same way the other sheet's header is styled:
# bold, thin box border, centered horizontally, top-aligned vertically.
$playerHeader = $playerSheet.Range("A1:D1")
$playerHeader.Font.Bold = $true
$playerHeader.Borders.LineStyle = 1
$playerHeader.HorizontalAlignment = -4108
$playerHeader.VerticalAlignment = -4160

# Data row - keep the ID as text (matches inlineStr "5924" in source data)
$idCell = $playerSheet.Cells.Item(2, 1)
$idCell.NumberFormat = "@"
$idCell.Value = "5924"
$playerSheet.Cells.Item(2, 2).Value = "Devon Philip Conway"
$playerSheet.Cells.Item(2, 3).Value = "Left Handed"
$playerSheet.Cells.Item(2, 4).Value = "Right Arm Medium"

# --- 2. Update the "ODI Batting" sheet: MATCH_CARD_LINK -> MATCH_CODE ---
# (Re-fetch the sheet reference by name now, AFTER the insert above, since
#  sheet handles are positional and would otherwise point at the new sheet.)
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$battingSheet.Range("D1").Value = "MATCH_CODE"

$matchCodes = @{
    2  = "4452"
    3  = "4453"
    4  = "4455"
    5  = "4636"
    6  = "4639"
    7  = "4642"
    8  = "4647"
    9  = "4648"
    10 = "4649"
    11 = "4669"
    12 = "4673"
    13 = "4676"
    14 = "4686"
    15 = "4688"
    16 = "4690"
    17 = "4692"
    18 = "4695"
    19 = "4697"
}

foreach ($row in $matchCodes.Keys) {
    $cell = $battingSheet.Cells.Item($row, 4)
    $cell.NumberFormat = "@"
    $cell.Value = $matchCodes[$row]
}

Write-Host "Applied Player Info sheet + MATCH_CODE column update"
